# Auto-generated edit script: updates market-price derived columns (H-N)
# across rows in multiple sheets, per the scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 374.55554
$ws.Range("I33").Value = 188.85715
$ws.Range("K33").Value = 188.85715
$ws.Range("M33").Value = 40.14285000000001
$ws.Range("H100").Value = 7565.1875
$ws.Range("I100").Value = 7399.3335
$ws.Range("K100").Value = 7399.3335
$ws.Range("M100").Value = -6858.3335
$ws.Range("H130").Value = 75000
$ws.Range("J130").Value = 75000
$ws.Range("L130").Value = 75000
$ws.Range("N130").Value = -85040

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 73973
$ws.Range("J28").Value = 200519
$ws.Range("L28").Value = 200519
$ws.Range("N28").Value = -200903
$ws.Range("H32").Value = 4931.6494
$ws.Range("I32").Value = 5051.877
$ws.Range("K32").Value = 5051.877
$ws.Range("M32").Value = -4764.877
$ws.Range("H61").Value = 2138.611
$ws.Range("I61").Value = 1823.5714
$ws.Range("J61").Value = 3241.25
$ws.Range("K61").Value = 1823.5714
$ws.Range("L61").Value = 3241.25
$ws.Range("M61").Value = -1611.5714
$ws.Range("N61").Value = -3665.25
$ws.Range("H99").Value = 73973
$ws.Range("J99").Value = 200519
$ws.Range("L99").Value = 200519
$ws.Range("N99").Value = -206509
$ws.Range("H122").Value = 5084.6787
$ws.Range("I122").Value = 4663.7144
$ws.Range("K122").Value = 13991.1432
$ws.Range("M122").Value = -11541.1432
$ws.Range("H136").Value = 2138.611
$ws.Range("I136").Value = 1823.5714
$ws.Range("J136").Value = 3241.25
$ws.Range("K136").Value = 5470.7142
$ws.Range("L136").Value = 9723.75
$ws.Range("M136").Value = -2920.7142
$ws.Range("N136").Value = -14823.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H50").Value = 70000
$ws.Range("J50").Value = 70000
$ws.Range("L50").Value = 70000
$ws.Range("N50").Value = -71148
$ws.Range("H94").Value = 7800
$ws.Range("I94").Value = 6000
$ws.Range("K94").Value = 6000
$ws.Range("M94").Value = -5549
$ws.Range("H99").Value = 38626
$ws.Range("I99").Value = 42363.6
$ws.Range("J99").Value = 1250
$ws.Range("K99").Value = 42363.6
$ws.Range("L99").Value = 1250
$ws.Range("M99").Value = -40865.6
$ws.Range("N99").Value = -4246
$ws.Range("H105").Value = 1993.25
$ws.Range("I105").Value = 1847.8
$ws.Range("K105").Value = 1847.8
$ws.Range("M105").Value = -100.8
$ws.Range("H107").Value = 3665.139
$ws.Range("I107").Value = 3410
$ws.Range("K107").Value = 3410
$ws.Range("M107").Value = -1490

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3261.2942
$ws.Range("I16").Value = 3079.6667
$ws.Range("K16").Value = 3079.6667
$ws.Range("M16").Value = -2792.6667
$ws.Range("H31").Value = 1547.4642
$ws.Range("I31").Value = 1459.5834
$ws.Range("K31").Value = 1459.5834
$ws.Range("M31").Value = -1164.5834
$ws.Range("H34").Value = 1547.4642
$ws.Range("I34").Value = 1459.5834
$ws.Range("K34").Value = 1459.5834
$ws.Range("M34").Value = -1257.5834
$ws.Range("H58").Value = 1527.5
$ws.Range("I58").Value = 884.4545000000001
$ws.Range("K58").Value = 884.4545000000001
$ws.Range("M58").Value = -681.4545000000001
$ws.Range("H99").Value = 3603.4285
$ws.Range("I99").Value = 4166.125
$ws.Range("K99").Value = 4166.125
$ws.Range("M99").Value = -2668.125
$ws.Range("H113").Value = 3261.2942
$ws.Range("I113").Value = 3079.6667
$ws.Range("K113").Value = 3079.6667
$ws.Range("M113").Value = -909.6667000000002
$ws.Range("H126").Value = 3603.4285
$ws.Range("I126").Value = 4166.125
$ws.Range("K126").Value = 12498.375
$ws.Range("M126").Value = -10028.375
$ws.Range("H134").Value = 3978.9666
$ws.Range("I134").Value = 3963.1538
$ws.Range("K134").Value = 11889.4614
$ws.Range("M134").Value = -9354.4614
$ws.Range("H136").Value = 1527.5
$ws.Range("I136").Value = 884.4545000000001
$ws.Range("K136").Value = 2653.3635
$ws.Range("M136").Value = -103.3635000000004

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 799.03125
$ws.Range("J5").Value = 1133.25
$ws.Range("L5").Value = 3399.75
$ws.Range("N5").Value = -3623.75
$ws.Range("H68").Value = 2161.077
$ws.Range("I68").Value = 942
$ws.Range("K68").Value = 2826
$ws.Range("M68").Value = -2015
$ws.Range("H71").Value = 2161.077
$ws.Range("I71").Value = 942
$ws.Range("K71").Value = 8478
$ws.Range("M71").Value = -4422
$ws.Range("H95").Value = 4995
$ws.Range("I95").Value = 4990
$ws.Range("K95").Value = 14970
$ws.Range("M95").Value = -12911
$ws.Range("H129").Value = 1962.6666
$ws.Range("I129").Value = 288
$ws.Range("J129").Value = 2800
$ws.Range("K129").Value = 864
$ws.Range("L129").Value = 8400
$ws.Range("M129").Value = 4136
$ws.Range("N129").Value = -18400
$ws.Range("H135").Value = 799.03125
$ws.Range("J135").Value = 1133.25
$ws.Range("L135").Value = 10199.25
$ws.Range("N135").Value = -15269.25
$ws.Range("H140").Value = 13893399
$ws.Range("J140").Value = 6999.857
$ws.Range("L140").Value = 20999.571
$ws.Range("N140").Value = -31359.571

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 1404250.8
$ws.Range("I11").Value = 1685001
$ws.Range("J11").Value = 562000
$ws.Range("K11").Value = 1685001
$ws.Range("L11").Value = 562000
$ws.Range("M11").Value = -1684862
$ws.Range("N11").Value = -562278
$ws.Range("H35").Value = 17500
$ws.Range("I35").Value = 17500
$ws.Range("K35").Value = 17500
$ws.Range("M35").Value = -17202
$ws.Range("H113").Value = 6256.2144
$ws.Range("I113").Value = 2945.25
$ws.Range("K113").Value = 2945.25
$ws.Range("M113").Value = -775.25

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 12255
$ws.Range("I30").Value = 12255
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 12255
$ws.Range("L30").Value = 0
$ws.Range("M30").ClearContents()
$ws.Range("N30").Value = -12147
$ws.Range("H61").Value = 2385.75
$ws.Range("I61").Value = 2431.1667
$ws.Range("K61").Value = 2431.1667
$ws.Range("M61").Value = -2229.1667
$ws.Range("H93").Value = 5998.8
$ws.Range("I93").Value = 6248.5
$ws.Range("J93").Value = 5000
$ws.Range("K93").Value = 6248.5
$ws.Range("L93").Value = 5000
$ws.Range("M93").Value = -5000.5
$ws.Range("N93").Value = -7496
$ws.Range("H100").Value = 374367
$ws.Range("I100").Value = 558050.5
$ws.Range("K100").Value = 558050.5
$ws.Range("M100").Value = -557509.5
$ws.Range("H113").Value = 2385.75
$ws.Range("I113").Value = 2431.1667
$ws.Range("K113").Value = 2431.1667
$ws.Range("M113").Value = -261.1667000000002

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1128.1428
$ws.Range("I113").Value = 699.6667
$ws.Range("J113").Value = 1449.5
$ws.Range("K113").Value = 2099.0001
$ws.Range("L113").Value = 4348.5
$ws.Range("M113").Value = 70.9998999999998
$ws.Range("N113").Value = -8688.5
$ws.Range("H114").Value = 39498
$ws.Range("J114").Value = 39498
$ws.Range("L114").Value = 39498
$ws.Range("N114").Value = -48176
$ws.Range("H126").Value = 7593.875
$ws.Range("I126").Value = 8407.286
$ws.Range("K126").Value = 25221.858
$ws.Range("M126").Value = -22751.858
$ws.Range("H132").Value = 1130.8438
$ws.Range("I132").Value = 730.7586
$ws.Range("K132").Value = 2192.2758
$ws.Range("M132").Value = 337.7242000000001
